# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'19.902.31"
$ws.Range("E2").Value = '  -8.35%  '

$ws.Range("D3").Value = "'1.403.62"
$ws.Range("E3").Value = '  -8.75%  '

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = '  +0.22%  '

$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").Value = "'272.63"
$ws.Range("E6").Value = '  -5.92%  '

$ws.Range("D7").Value = "'0.3674"
$ws.Range("E7").Value = '  -6.48%  '

$ws.Range("D8").Value = "'0.3109"
$ws.Range("E8").Value = '  -2.74%  '

$ws.Range("D9").Value = "'39.46"
$ws.Range("E9").Value = '  -8.88%  '

$ws.Range("D10").Value = "'1.006"
$ws.Range("E10").Value = '  -6.26%  '

$ws.Range("D11").Value = "'0.06468"
$ws.Range("E11").Value = '  -10.08%  '

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").Value = "'5.413"
$ws.Range("E13").Value = '  -5.94%  '

$ws.Range("D14").Value = "'17.44"
$ws.Range("E14").Value = '  -5.82%  '

$ws.Range("D15").Value = "'6.146"
$ws.Range("E15").Value = '  -7.43%  '

$ws.Range("D16").Value = "'1.402.66"
$ws.Range("E16").Value = '  -5.93%  '

$ws.Range("D17").Value = "'0.00001006"
$ws.Range("E17").Value = '  -8.36%  '

$ws.Range("D18").Value = "'0.05667"
$ws.Range("E18").Value = '  -14.29%  '

$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = '  +0.23%  '

$ws.Range("D20").Value = "'70.31"
$ws.Range("E20").Value = '  -16.64%  '

$ws.Range("D21").Value = "'5.557"
$ws.Range("E21").Value = '  -9.76%  '

$ws.Range("D22").Value = "'14.67"
$ws.Range("E22").Value = '  -5.72%  '

$ws.Range("D23").Value = "'10.99"
$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("D24").Value = "'2.278"
$ws.Range("E24").Value = '  -3.61%  '

$ws.Range("D25").Value = "'19.927.44"
$ws.Range("E25").Value = '  -8.20%  '

$ws.Range("D26").Value = "'2.232"
$ws.Range("E26").Value = '  -6.68%  '

$ws.Range("D27").Value = "'135.19"
$ws.Range("E27").Value = '  -11.02%  '

$ws.Range("D28").Value = "'16.82"
$ws.Range("E28").Value = '  -9.14%  '

$ws.Range("D29").Value = "'1.560.40"
$ws.Range("E29").Value = '  -6.41%  '

$ws.Range("D30").Value = "'108.84"
$ws.Range("E30").Value = '  -7.62%  '

$ws.Range("D31").Value = "'4.084"
$ws.Range("E31").Value = '  -16.02%  '

$ws.Range("D32").Value = "'5.285"
$ws.Range("E32").Value = '  -13.29%  '

$ws.Range("D33").Value = "'0.8091"
$ws.Range("E33").Value = '  -16.75%  '

$ws.Range("D34").Value = "'0.07644"
$ws.Range("E34").Value = '  -5.87%  '

$ws.Range("D35").Value = "'8.366"
$ws.Range("E35").Value = '  -2.74%  '

$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = "'1.440"
$ws.Range("E36").Value = '  -3.81%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = "'0.05795"
$ws.Range("E37").Value = '  -3.35%  '

$ws.Range("D38").Value = "'4.797"
$ws.Range("E38").Value = '  -8.04%  '

$ws.Range("D39").Value = "'1.001"
$ws.Range("E39").Value = '  +0.16%  '

$ws.Range("D40").Value = "'0.02062"
$ws.Range("E40").Value = '  -8.02%  '

$ws.Range("D41").Value = "'0.1900"
$ws.Range("E41").Value = '  -7.45%  '

$ws.Range("D42").Value = "'10.34"
$ws.Range("E42").Value = '  -8.84%  '

$ws.Range("D43").Value = "'1.086"
$ws.Range("E43").Value = '  -8.35%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = "'12.25"
$ws.Range("E44").Value = '  -7.98%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = "'0.5255"
$ws.Range("E45").Value = '  -9.94%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = "'3.500"
$ws.Range("E46").Value = '  -6.14%  '

$ws.Range("D47").Value = "'0.5093"
$ws.Range("E47").Value = '  -8.92%  '

$ws.Range("D48").Value = "'111.92"
$ws.Range("E48").Value = '  -3.78%  '

$ws.Range("D49").Value = "'1.757"
$ws.Range("E49").Value = '  -7.33%  '

$ws.Range("D50").Value = "'1.033"
$ws.Range("E50").Value = '  -11.17%  '

$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = '  +0.24%  '
